# Weekly refresh of the Tomate (Macroferia Regional de Talca) price series.
# Rows 460-555 hold the historical weekly records; this update:
#   1) shifts all of those records down by 3 rows (460-555 -> 463-558), and
#   2) inserts 3 brand-new records (the latest week) at rows 460-462.
# Columns D,I,J,K,L,M,N,O,P,Q carry the per-record data; A,B,C,E,F,G,H,R are
# constant across the whole block and are left untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$firstRow = 460
$lastRow = 555
$shift = 3

# The leading/trailing columns (Mercado ID, Mercado, Region, Codreg,
# Categoria ID, Categoria, Variedad, Clasificacion) are constant for the
# whole block, so grab them once from an existing row to stamp the 3
# brand-new rows (556-558) that don't exist yet.
$constA = $ws.Cells.Item($firstRow, 1).Value2
$constB = $ws.Cells.Item($firstRow, 2).Value2
$constC = $ws.Cells.Item($firstRow, 3).Value2
$constE = $ws.Cells.Item($firstRow, 5).Value2
$constF = $ws.Cells.Item($firstRow, 6).Value2
$constG = $ws.Cells.Item($firstRow, 7).Value2
$constH = $ws.Cells.Item($firstRow, 8).Value2
$constR = $ws.Cells.Item($firstRow, 18).Value2

# --- 1) Snapshot the existing records (460-555) before overwriting anything ---
$records = @()
for ($r = $firstRow; $r -le $lastRow; $r++) {
    $rec = @{
        D = $ws.Cells.Item($r, 4).Value2
        I = $ws.Cells.Item($r, 9).Value2
        J = $ws.Cells.Item($r, 10).Value2
        K = $ws.Cells.Item($r, 11).Value2
        L = $ws.Cells.Item($r, 12).Value2
        M = $ws.Cells.Item($r, 13).Value2
        N = $ws.Cells.Item($r, 14).Value2
        O = $ws.Cells.Item($r, 15).Value2
        P = $ws.Cells.Item($r, 16).Value2
        Q = $ws.Cells.Item($r, 17).Value2
    }
    $records += $rec
}

# --- 2) Write the snapshot back out shifted down by 3 rows (460-555 -> 463-558) ---
for ($i = 0; $i -lt $records.Count; $i++) {
    $rec = $records[$i]
    $r = $firstRow + $shift + $i

    $ws.Cells.Item($r, 1).Value = $constA
    $ws.Cells.Item($r, 2).Value = $constB
    $ws.Cells.Item($r, 3).Value = $constC
    $ws.Cells.Item($r, 4).Value = $rec.D
    $ws.Cells.Item($r, 5).Value = $constE
    $ws.Cells.Item($r, 6).Value = $constF
    $ws.Cells.Item($r, 7).Value = $constG
    $ws.Cells.Item($r, 8).Value = $constH
    $ws.Cells.Item($r, 9).Value = $rec.I
    $ws.Cells.Item($r, 10).Value = $rec.J
    $ws.Cells.Item($r, 11).Value = $rec.K
    $ws.Cells.Item($r, 12).Value = $rec.L
    $ws.Cells.Item($r, 13).Value = $rec.M
    $ws.Cells.Item($r, 14).Value = $rec.N
    $ws.Cells.Item($r, 15).Value = $rec.O
    $ws.Cells.Item($r, 16).Value = $rec.P
    $ws.Cells.Item($r, 17).Value = $rec.Q
    $ws.Cells.Item($r, 18).Value = $constR
}

# --- 3) Fill the 3 newly-opened rows (460-462) with this week's new records ---
$newRecords = @(
    @{ D = 44637; J = 2500; K = 18000; L = 18000; M = 18000; N = "$/bandeja 18 kilos"; O = "Región de Arica y Parinacota"; P = 1000; Q = 18 },
    @{ D = 44637; J = 2000; K = 15000; L = 15000; M = 15000; N = "$/bandeja 18 kilos"; O = "Región del Maule";             P = 833;  Q = 18 },
    @{ D = 44637; J = 2000; K = 8000;  L = 8000;  M = 8000;  N = "$/caja 15 kilos";    O = "Región del Maule";             P = 533;  Q = 15 }
)

for ($i = 0; $i -lt $newRecords.Count; $i++) {
    $rec = $newRecords[$i]
    $r = $firstRow + $i

    $ws.Cells.Item($r, 1).Value = $constA
    $ws.Cells.Item($r, 2).Value = $constB
    $ws.Cells.Item($r, 3).Value = $constC
    $ws.Cells.Item($r, 4).Value = $rec.D
    $ws.Cells.Item($r, 5).Value = $constE
    $ws.Cells.Item($r, 6).Value = $constF
    $ws.Cells.Item($r, 7).Value = $constG
    $ws.Cells.Item($r, 8).Value = $constH
    $ws.Cells.Item($r, 9).Value = "Primera"
    $ws.Cells.Item($r, 10).Value = $rec.J
    $ws.Cells.Item($r, 11).Value = $rec.K
    $ws.Cells.Item($r, 12).Value = $rec.L
    $ws.Cells.Item($r, 13).Value = $rec.M
    $ws.Cells.Item($r, 14).Value = $rec.N
    $ws.Cells.Item($r, 15).Value = $rec.O
    $ws.Cells.Item($r, 16).Value = $rec.P
    $ws.Cells.Item($r, 17).Value = $rec.Q
    $ws.Cells.Item($r, 18).Value = $constR
}
